# Modifications pour utiliser XGBClassifier et ajuster les prédictions
#
# Sheet "Valeurs réelles" (sheet1): rename the _S+1/_S+2/_S+3 headers to
# *_class, and replace the numeric forecast columns (C:E) with integer
# class labels (XGBClassifier output) instead of the old float prices.
#
# Sheet "Prédictions" (sheet2): zero out the predicted columns (B:D) —
# the classifier no longer predicts a continuous price — except for the
# one row where the model predicts the negative class (-2).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Valeurs réelles")
$ws2 = $wb.Worksheets.Item("Prédictions")

# ---------------------------------------------------------------------
# Sheet 1 - "Valeurs réelles": header renames
# ---------------------------------------------------------------------
$ws1.Range("C1").Value = "PRIX EXP POMME BRAEBURN FRANCE 170/220G CAT.I PLATEAU 1RG_S+1_class"
$ws1.Range("D1").Value = "PRIX EXP POMME BRAEBURN FRANCE 170/220G CAT.I PLATEAU 1RG_S+2_class"
$ws1.Range("E1").Value = "PRIX EXP POMME BRAEBURN FRANCE 170/220G CAT.I PLATEAU 1RG_S+3_class"

# ---------------------------------------------------------------------
# Sheet 1 - "Valeurs réelles": new class values for columns C, D, E
# (row -> value) for rows 2..28
# ---------------------------------------------------------------------
$class_C = @{2=4; 3=2; 4=2; 5=2; 6=2; 7=2; 8=2; 9=2; 10=2; 11=2; 12=2; 13=2; 14=2; 15=2; 16=2; 17=2; 18=2; 19=2; 20=2; 21=2; 22=4; 23=2; 24=2; 25=2; 26=2; 27=4; 28=2}
$class_D = @{2=2; 3=2; 4=2; 5=2; 6=2; 7=2; 8=2; 9=2; 10=2; 11=2; 12=2; 13=2; 14=2; 15=2; 16=2; 17=2; 18=2; 19=2; 20=2; 21=4; 22=2; 23=2; 24=2; 25=2; 26=4; 27=2; 28=2}
$class_E = @{2=2; 3=2; 4=2; 5=2; 6=2; 7=2; 8=2; 9=2; 10=2; 11=2; 12=2; 13=2; 14=2; 15=2; 16=2; 17=2; 18=2; 19=2; 20=4; 21=2; 22=2; 23=2; 24=2; 25=4; 26=2; 27=2; 28=2}

foreach ($r in 2..28) {
    $ws1.Cells.Item($r, 3).Value = $class_C[$r]
    $ws1.Cells.Item($r, 4).Value = $class_D[$r]
    $ws1.Cells.Item($r, 5).Value = $class_E[$r]
}

# ---------------------------------------------------------------------
# Sheet 2 - "Prédictions": predicted columns B, C, D all become 0,
# except row 22 (DATE_INTERROGATION = 45642) where PRED_S1 (col B) is -2.
# ---------------------------------------------------------------------
foreach ($r in 2..28) {
    $ws2.Cells.Item($r, 2).Value = 0
    $ws2.Cells.Item($r, 3).Value = 0
    $ws2.Cells.Item($r, 4).Value = 0
}
$ws2.Cells.Item(22, 2).Value = -2
